$wb = $excel.ActiveWorkbook

# "Latest HO Xliff Generate Date" on Overview and "Correspond Handoff
# Datetime" on de-de both shared the same timestamp string
# (2016-09-07 03:22:23); the report refresh bumped that shared value to
# 2016-09-07 03:23:13, so both cells need to move together.
$wsOverview = $wb.Worksheets.Item("Overview")
$wsOverview.Range("G2").Value = "2016-09-07 03:23:13"

$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsZhCn.Range("H2").Value = "2016-09-07 03:23:03"
$wsZhCn.Range("K2").Value = "2016-09-07 03:23:34"

$wsDeDe = $wb.Worksheets.Item("de-de")
$wsDeDe.Range("H2").Value = "2016-09-07 03:23:13"
$wsDeDe.Range("K2").Value = "2016-09-07 03:23:42"
